# Updated cryptos list on Sun Jun 23 13:43:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text (they hold values
# like "588.55" or "1.00" that Excel would otherwise auto-coerce into
# numbers), then restore the default "Normal" style afterwards so no
# stray number-format/style ends up attached to the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.310.91"
$ws.Range("E2").Value = "  +0.04%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.498.76"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "588.55"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6 - Solana
$ws.Range("D6").Value = "133.68"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "7.67"
$ws.Range("E9").Value = "  +6.37%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.27%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.387"
$ws.Range("E11").Value = "  -0.28%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.094.35"
$ws.Range("E12").Value = "  +0.15%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.18%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -0.85%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.497.83"
$ws.Range("E15").Value = "  +0.10%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "64.256.05"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "24.86"
$ws.Range("E17").Value = "  -3.90%  "

# Row 18 - Uniswap
$ws.Range("E18").Value = "  +0.41%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -0.59%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -1.91%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "386.24"
$ws.Range("E21").Value = "  -0.33%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  +1.83%  "

# Row 23 - WrappedeETH
$ws.Range("D23").Value = "3.637.55"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "74.35"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.03%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  +0.17%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +1.64%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.03%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  -2.11%  "

# Rows 30 & 31 swap positions: Fetch.AI <-> PancakeSwap.
# Row 30 becomes PancakeSwap, row 31 becomes Fetch.AI (with updated data).
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.24"
$ws.Range("E30").Value = "  +0.53%  "

$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.49"
$ws.Range("E31").Value = "  +0.67%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -1.80%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  +3.92%  "

# Row 34 - RenzoRestakedETH
$ws.Range("D34").Value = "3.526.45"
$ws.Range("E34").Value = "  +0.30%  "

# Row 35 - USDe
$ws.Range("E35").Value = "  -0.02%  "

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "23.26"
$ws.Range("E36").Value = "  -1.22%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +3.56%  "

# Row 38 - Aptos
$ws.Range("E38").Value = "  +0.00%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -1.46%  "

# Row 40 - Monero
$ws.Range("D40").Value = "164.47"
$ws.Range("E40").Value = "  +0.53%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.0783"
$ws.Range("E41").Value = "  +0.02%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "0.807"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.03%  "

# Row 44 - Filecoin
$ws.Range("E44").Value = "  -0.77%  "

# Row 45 - ONDO
$ws.Range("E45").Value = "  +0.89%  "

# Row 46 - EnergySwap
$ws.Range("D46").Value = "24.35"
$ws.Range("E46").Value = "  -6.06%  "

# Row 47 - Stacks
$ws.Range("D47").Value = "1.64"
$ws.Range("E47").Value = "  -1.13%  "

# Row 48 - Maker
$ws.Range("D48").Value = "2.427.35"
$ws.Range("E48").Value = "  -2.25%  "

# Row 49 - SuiNetwork
$ws.Range("E49").Value = "  +2.52%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  +0.07%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  -1.55%  "

# Restore the default style on the whole range so we don't leave a
# lingering text-number-format/quote-prefix style on any cell.
$dataRange.Style = "Normal"
